# Auto-generated Excel COM-interop script
# Applies the numeric corrections to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as captured by the scheduled runner diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 38280
$ws.Range("J133").Value = 38280
$ws.Range("L133").Value = 38280
$ws.Range("N133").Value = -48400
$ws.Range("H134").Value = 29166.666
$ws.Range("J134").Value = 29166.666
$ws.Range("L134").Value = 29166.666
$ws.Range("N134").Value = -39306.666
$ws.Range("H136").Value = 36468
$ws.Range("J136").Value = 36468
$ws.Range("L136").Value = 36468
$ws.Range("N136").Value = -46668
$ws.Range("H139").Value = 35533.9
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 35533.9
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 35533.9
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -45813.9
$ws.Range("H140").Value = 33853.332
$ws.Range("J140").Value = 33853.332
$ws.Range("L140").Value = 33853.332
$ws.Range("N140").Value = -44213.332
$ws.Range("H141").Value = 1502354.5
$ws.Range("I141").Value = 1580899.5
$ws.Range("K141").Value = 4742698.5
$ws.Range("M141").Value = -4737518.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3422.88
$ws.Range("I32").Value = 3229.7732
$ws.Range("J32").Value = 9666.666999999999
$ws.Range("K32").Value = 3229.7732
$ws.Range("L32").Value = 9666.666999999999
$ws.Range("M32").Value = -2942.7732
$ws.Range("N32").Value = -10240.667
$ws.Range("H74").Value = 143972
$ws.Range("I74").Value = 167837.33
$ws.Range("K74").Value = 167837.33
$ws.Range("M74").Value = -166963.33
$ws.Range("H77").Value = 143972
$ws.Range("I77").Value = 167837.33
$ws.Range("K77").Value = 839186.6499999999
$ws.Range("M77").Value = -834818.6499999999
$ws.Range("H102").Value = 1087.4286
$ws.Range("I102").Value = 1087.4286
$ws.Range("K102").Value = 1087.4286
$ws.Range("M102").Value = 534.5714
$ws.Range("H134").Value = 31964.5
$ws.Range("J134").Value = 31964.5
$ws.Range("L134").Value = 31964.5
$ws.Range("N134").Value = -42104.5
$ws.Range("H135").Value = 28104.143
$ws.Range("J135").Value = 28104.143
$ws.Range("L135").Value = 28104.143
$ws.Range("N135").Value = -38244.143
$ws.Range("H139").Value = 80984.44500000001
$ws.Range("J139").Value = 80984.44500000001
$ws.Range("L139").Value = 80984.44500000001
$ws.Range("N139").Value = -91264.44500000001
$ws.Range("H140").Value = 42971.6
$ws.Range("J140").Value = 42971.6
$ws.Range("L140").Value = 42971.6
$ws.Range("N140").Value = -53331.6
$ws.Range("H141").Value = 45323.375
$ws.Range("J141").Value = 45323.375
$ws.Range("L141").Value = 45323.375
$ws.Range("N141").Value = -55683.375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 54029.668
$ws.Range("J81").Value = 54029.668
$ws.Range("L81").Value = 54029.668
$ws.Range("N81").Value = -56151.668
$ws.Range("H84").Value = 54029.668
$ws.Range("J84").Value = 54029.668
$ws.Range("L84").Value = 162089.004
$ws.Range("N84").Value = -172697.004
$ws.Range("H132").Value = 34000
$ws.Range("J132").Value = 34000
$ws.Range("L132").Value = 34000
$ws.Range("N132").Value = -44120
$ws.Range("H135").Value = 30727.182
$ws.Range("J135").Value = 30727.182
$ws.Range("L135").Value = 30727.182
$ws.Range("N135").Value = -40867.182
$ws.Range("H137").Value = 43646.188
$ws.Range("J137").Value = 43646.188
$ws.Range("L137").Value = 43646.188
$ws.Range("N137").Value = -53846.188
$ws.Range("H138").Value = 37888.777
$ws.Range("J138").Value = 39874.875
$ws.Range("L138").Value = 39874.875
$ws.Range("N138").Value = -50154.875

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2250.6667
$ws.Range("I31").Value = 1366.6786
$ws.Range("J31").Value = 3706.647
$ws.Range("K31").Value = 1366.6786
$ws.Range("L31").Value = 3706.647
$ws.Range("M31").Value = -1071.6786
$ws.Range("N31").Value = -4296.647
$ws.Range("H34").Value = 2250.6667
$ws.Range("I34").Value = 1366.6786
$ws.Range("J34").Value = 3706.647
$ws.Range("K34").Value = 1366.6786
$ws.Range("L34").Value = 3706.647
$ws.Range("M34").Value = -1164.6786
$ws.Range("N34").Value = -4110.647
$ws.Range("H135").Value = 32616.363
$ws.Range("J135").Value = 32616.363
$ws.Range("L135").Value = 32616.363
$ws.Range("N135").Value = -42756.363
$ws.Range("H138").Value = 37576.668
$ws.Range("J138").Value = 37576.668
$ws.Range("L138").Value = 37576.668
$ws.Range("N138").Value = -47856.668
$ws.Range("H141").Value = 44821.09
$ws.Range("J141").Value = 44821.09
$ws.Range("L141").Value = 44821.09
$ws.Range("N141").Value = -55181.09

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5777.778
$ws.Range("I3").Value = 3750
$ws.Range("J3").Value = 7400
$ws.Range("K3").Value = 11250
$ws.Range("L3").Value = 22200
$ws.Range("M3").Value = -11138
$ws.Range("N3").Value = -22424

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120
$ws.Range("H134").Value = 22119
$ws.Range("J134").Value = 22119
$ws.Range("L134").Value = 66357
$ws.Range("N134").Value = -71427
$ws.Range("H135").Value = 30793.086
$ws.Range("J135").Value = 30793.086
$ws.Range("L135").Value = 30793.086
$ws.Range("N135").Value = -40933.086
$ws.Range("H140").Value = 40780
$ws.Range("J140").Value = 40780
$ws.Range("L140").Value = 40780
$ws.Range("N140").Value = -51140
$ws.Range("H141").Value = 49117.582
$ws.Range("J141").Value = 49117.582
$ws.Range("L141").Value = 49117.582
$ws.Range("N141").Value = -59477.582

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 29928.625
$ws.Range("J134").Value = 29928.625
$ws.Range("L134").Value = 29928.625
$ws.Range("N134").Value = -40068.625
$ws.Range("H138").Value = 27929
$ws.Range("J138").Value = 27929
$ws.Range("L138").Value = 27929
$ws.Range("N138").Value = -38209

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 30998.166
$ws.Range("J46").Value = 30998.166
$ws.Range("L46").Value = 30998.166
$ws.Range("N46").Value = -31460.166
$ws.Range("H133").Value = 28000
$ws.Range("J133").Value = 28000
$ws.Range("L133").Value = 28000
$ws.Range("N133").Value = -38120
$ws.Range("H134").Value = 30998.166
$ws.Range("J134").Value = 30998.166
$ws.Range("L134").Value = 92994.49800000001
$ws.Range("N134").Value = -98064.49800000001
$ws.Range("H135").Value = 31893
$ws.Range("J135").Value = 31893
$ws.Range("L135").Value = 31893
$ws.Range("N135").Value = -42033
$ws.Range("H137").Value = 38629
$ws.Range("J137").Value = 38629
$ws.Range("L137").Value = 38629
$ws.Range("N137").Value = -48829
$ws.Range("H140").Value = 30228.572
$ws.Range("J140").Value = 30228.572
$ws.Range("L140").Value = 30228.572
$ws.Range("N140").Value = -40588.572
